# [Kadastro App] Yeni kayit eklendi: 2959
# Adds a new record (Kayit No 2959, Erdemli birimi, 09-09-2025, 3B isi)
# both to the master "Kayitlar" sheet (next free row) and to the
# per-birim "Erdemli" filtered sheet (next free row).

$wb = $excel.ActiveWorkbook

function Add-KayitRow {
    param($Worksheet, $Row, $KayitNo, $Tarih, $Birim, $ParselSayisi, $Is, $Personeller)

    $rowRange = $Worksheet.Range($Worksheet.Cells.Item($Row, 1), $Worksheet.Cells.Item($Row, 6))
    # Keep the new row text-typed (matching every other data cell in this
    # workbook, which is stored as text even for numeric-looking values).
    $rowRange.NumberFormat = "@"

    $Worksheet.Cells.Item($Row, 1).Value = $KayitNo
    $Worksheet.Cells.Item($Row, 2).Value = $Tarih
    $Worksheet.Cells.Item($Row, 3).Value = $Birim
    $Worksheet.Cells.Item($Row, 4).Value = $ParselSayisi
    $Worksheet.Cells.Item($Row, 5).Value = $Is
    $Worksheet.Cells.Item($Row, 6).Value = $Personeller
}

# --- "Kayitlar" (master records) sheet: new row at 29 ---
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-KayitRow $wsKayitlar 29 "2959" "2025-09-09" "Erdemli" "" "3B" "ÖZKAN AKBAŞ (Mühendis), ENDER NUSRET ÖNAL GÜLSOY (Kontrol Memuru)"

# --- "Erdemli" (per-birim filtered) sheet: new row at 28 ---
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-KayitRow $wsErdemli 28 "2959" "2025-09-09" "Erdemli" "" "3B" "ÖZKAN AKBAŞ (Mühendis), ENDER NUSRET ÖNAL GÜLSOY (Kontrol Memuru)"
